$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 109 (the existing row 109 and everything
# below it shifts down by one) for the new "powder for solution for
# intraocular irrigation" dose-form entry.
$ws.Rows.Item(109).Insert()

$ws.Cells.Item(109, 2).Value = "powder for solution for intraocular injection"
$ws.Cells.Item(109, 3).Value = "Solid sterile preparation consisting of one or more powders, including freeze-dried powders, intended to be dissolved in the specified liquid to obtain a solution for intraocular irrigation."
$ws.Cells.Item(109, 5).Value = "Poudre pour solution pour irrigation intraoculaire"

# Match the formatting used elsewhere in the sheet for this kind of row:
# column A blank/highlighted, column B bold, column C wrapped text, column E
# as text, and no value/cell in column D.
$ws.Cells.Item(109, 1).Interior.ColorIndex = 6
$ws.Cells.Item(109, 2).Font.Bold = $true
$ws.Cells.Item(109, 3).WrapText = $true
$ws.Cells.Item(109, 5).NumberFormat = "@"
$ws.Cells.Item(109, 4).Clear()
$ws.Rows.Item(109).RowHeight = 28.8

# Keep the sheet's defined name / filter range in sync with the extra row.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Sheet1!`$E`$1:`$E`$151"

$ws.Range("A109").Select() | Out-Null
